# Swap the step/expected-result content between the TC3 and TC4 blocks.
# TC3 (row 25) currently holds the "analisar prestação de contas" content.
# TC4 (row 32) currently holds the "cancelar diária" content.
# After the edit, TC3 should hold the "cancelar diária" content and
# TC4 should hold the "analisar prestação de contas" content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Steps = $ws.Range("B25").Value()
$tc3Expected = $ws.Range("D25").Value()
$tc4Steps = $ws.Range("B32").Value()
$tc4Expected = $ws.Range("D32").Value()

$ws.Range("B25").Value = $tc4Steps
$ws.Range("D25").Value = $tc4Expected

$ws.Range("B32").Value = $tc3Steps
$ws.Range("D32").Value = $tc3Expected
